# "classe Simulation en cours"
#
# Fix the typo "crée" -> "créée" in the "Structure des classes" intro
# paragraph and append a new sentence about the initial on-grid
# placement of each character. Word also leaves its usual collapsed
# "_GoBack" bookmark at the point of the last edit, so we reproduce
# that too.

$d = $word.ActiveDocument

$old = "Une classe pour chaque type personnage a été crée. "
# Append the fixed sentence plus the new one, with a throw-away
# trailing marker character ("~") so the insertion point where we will
# drop the bookmark isn't sitting exactly on the paragraph-end boundary
# while we create it.
$newWithMarker = "Une classe pour chaque type personnage a été créée.  De plus, l’emplacement initial de chaque personnage est inscrit sur la grille à l’aide des caractères respectifs.~"

$d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $newWithMarker, 2)

# Locate the marker and collapse a range to just before it - that is
# exactly where the paragraph will end once the marker is removed.
$markerRange = $d.Content
$markerRange.Find.Execute("~", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$markerRange.Collapse(1)

$d.Bookmarks.Add("_GoBack", $markerRange)

# Remove the throw-away marker character now that the bookmark is anchored.
$markerRange2 = $d.Content
$markerRange2.Find.Execute("~", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$markerRange2.Text = ""
